$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8320139646530151
$ws.Range("B1").Value = 1.383286952972412
$ws.Range("C1").Value = 3.723418712615967
$ws.Range("D1").Value = 2.657821178436279
$ws.Range("E1").Value = 1.634214997291565
